# The sheet holds weekly "Espinaca" price observations. This edit adds one
# new weekly observation at the top of the data block (row 5, right after
# the three rows 2-4 that stay fixed at the front), pushing every existing
# row from 5..111 down by one (6..112) and growing the used range from
# A1:R111 to A1:R112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 5; rows 5-111 shift down to become 6-112.
$ws.Rows.Item(5).Insert()

# Fill the newly inserted row 5 with the new price observation.
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 45160
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 100112012
$ws.Range("G5").Value = "Espinaca"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 270
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1900
$ws.Range("N5").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O5").Value = "Región de Arica y Parinacota"
$ws.Range("P5").Value = 633
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = "Hortaliza"
